$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell B93 with the new trim description
$ws.Range("B93").Value = "RX 350 AWD F SPORT BLACK LINE SPECIAL EDITION"

# Widen column B to fit the longer text (stored XML width of 55 equals a
# ColumnWidth of 55 - 5/6, since Excel's COM ColumnWidth excludes the fixed
# 5px cell-padding baked into the persisted <col width=".."/> value)
$ws.Columns.Item(2).ColumnWidth = 325/6

# Update the visible window / selection to reflect where the user was working
$ws.Application.ActiveWindow.ScrollRow = 65
$ws.Range("B93").Select()
